# Apply weekly data roll: insert two new rows of "Espinaca" price data
# at row 114 (pushing the existing rows 114-215 down to 116-217).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 114; this shifts existing rows 114-215
# down to 116-217, matching the target layout exactly (no data loss).
$ws.Rows("114:115").Insert()

# --- New row 114 ---
$ws.Cells.Item(114, 1).Value = 9
$ws.Cells.Item(114, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(114, 3).Value = "Metropolitana"
$ws.Cells.Item(114, 4).Value = 44447
$ws.Cells.Item(114, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(114, 5).Value = 13
$ws.Cells.Item(114, 6).Value = 100112012
$ws.Cells.Item(114, 7).Value = "Espinaca"
$ws.Cells.Item(114, 8).Value = "Sin especificar"
$ws.Cells.Item(114, 9).Value = "Primera"
$ws.Cells.Item(114, 10).Value = 106
$ws.Cells.Item(114, 11).Value = 7000
$ws.Cells.Item(114, 12).Value = 8000
$ws.Cells.Item(114, 13).Value = 7500
$ws.Cells.Item(114, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(114, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(114, 16).Value = 750
$ws.Cells.Item(114, 17).Value = 10
$ws.Cells.Item(114, 18).Value = "Hortaliza"

# --- New row 115 ---
$ws.Cells.Item(115, 1).Value = 9
$ws.Cells.Item(115, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(115, 3).Value = "Metropolitana"
$ws.Cells.Item(115, 4).Value = 44447
$ws.Cells.Item(115, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(115, 5).Value = 13
$ws.Cells.Item(115, 6).Value = 100112012
$ws.Cells.Item(115, 7).Value = "Espinaca"
$ws.Cells.Item(115, 8).Value = "Sin especificar"
$ws.Cells.Item(115, 9).Value = "Segunda"
$ws.Cells.Item(115, 10).Value = 52
$ws.Cells.Item(115, 11).Value = 6000
$ws.Cells.Item(115, 12).Value = 6000
$ws.Cells.Item(115, 13).Value = 6000
$ws.Cells.Item(115, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(115, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(115, 16).Value = 600
$ws.Cells.Item(115, 17).Value = 10
$ws.Cells.Item(115, 18).Value = "Hortaliza"
